# Atualizacao de bases das ligas - 07-03-2024 23:43
# Updates match odds/results data: swaps rows 91/92, and refreshes rows 179-183
# (new fixture inserted, two old fixtures resolved with final scores & updated closing odds,
#  one stale fixture replaced, and a new upcoming fixture appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91
$ws.Range('A91').Value = 89
$ws.Range('B91').Value = 6924569
$ws.Range('C91').Value = 'Mexico Liga de Expansion'
$ws.Range('D91').Value = 'Mexico Liga de Expansion'
$ws.Range('E91').Value = 45214.92013888889
$ws.Range('F91').Value = 'Venados FC'
$ws.Range('G91').Value = 'Dorados'
$ws.Range('H91').Value = 4
$ws.Range('I91').Value = 1
$ws.Range('J91').Value = 'H'
$ws.Range('K91').Value = 1.615
$ws.Range('L91').Value = 4
$ws.Range('M91').Value = 4.5
$ws.Range('N91').Value = 1.5
$ws.Range('O91').Value = 4.75
$ws.Range('P91').Value = 5.75
$ws.Range('Q91').Value = -1.25
$ws.Range('R91').Value = 1.925
$ws.Range('S91').Value = 1.875
$ws.Range('T91').Value = 3
$ws.Range('U91').Value = 1.75
$ws.Range('V91').Value = 1.95
$ws.Range('W91').Value = 0.5
$ws.Range('X91').Value = -1
$ws.Range('Y91').Value = -1
$ws.Range('Z91').Value = 0.925
$ws.Range('AA91').Value = -1
$ws.Range('AB91').Value = 0.75
$ws.Range('AC91').Value = -1

# Row 92
$ws.Range('A92').Value = 90
$ws.Range('B92').Value = 6924568
$ws.Range('C92').Value = 'Mexico Liga de Expansion'
$ws.Range('D92').Value = 'Mexico Liga de Expansion'
$ws.Range('E92').Value = 45214.92013888889
$ws.Range('F92').Value = 'Atletico Morelia'
$ws.Range('G92').Value = 'Atlante'
$ws.Range('H92').Value = 0
$ws.Range('I92').Value = 1
$ws.Range('J92').Value = 'A'
$ws.Range('K92').Value = 2.4
$ws.Range('L92').Value = 3
$ws.Range('M92').Value = 2.875
$ws.Range('N92').Value = 2.7
$ws.Range('O92').Value = 3.1
$ws.Range('P92').Value = 2.8
$ws.Range('Q92').Value = 0
$ws.Range('R92').Value = 1.85
$ws.Range('S92').Value = 1.95
$ws.Range('T92').Value = 2.25
$ws.Range('U92').Value = 1.975
$ws.Range('V92').Value = 1.725
$ws.Range('W92').Value = -1
$ws.Range('X92').Value = -1
$ws.Range('Y92').Value = 1.8
$ws.Range('Z92').Value = -1
$ws.Range('AA92').Value = 0.95
$ws.Range('AB92').Value = -1
$ws.Range('AC92').Value = 0.7250000000000001

# Row 179
$ws.Range('A179').Value = 177
$ws.Range('B179').Value = 7641686
$ws.Range('C179').Value = 'Mexico Liga de Expansion'
$ws.Range('D179').Value = 'Mexico Liga de Expansion'
$ws.Range('E179').Value = 45353.83680555555
$ws.Range('F179').Value = 'Cancun FC'
$ws.Range('G179').Value = 'Club Atletico La Paz'
$ws.Range('H179').Value = 2
$ws.Range('I179').Value = 2
$ws.Range('J179').Value = 'D'
$ws.Range('K179').Value = 2
$ws.Range('L179').Value = 3.1
$ws.Range('M179').Value = 3.4
$ws.Range('N179').Value = 1.615
$ws.Range('O179').Value = 3.6
$ws.Range('P179').Value = 6
$ws.Range('Q179').Value = -0.75
$ws.Range('R179').Value = 1.825
$ws.Range('S179').Value = 1.975
$ws.Range('T179').Value = 2.5
$ws.Range('U179').Value = 2
$ws.Range('V179').Value = 1.8
$ws.Range('W179').Value = -1
$ws.Range('X179').Value = 2.6
$ws.Range('Y179').Value = -1
$ws.Range('Z179').Value = -1
$ws.Range('AA179').Value = 0.9750000000000001
$ws.Range('AB179').Value = 1
$ws.Range('AC179').Value = -1

# Row 180
$ws.Range('A180').Value = 178
$ws.Range('B180').Value = 7640647
$ws.Range('C180').Value = 'Mexico Liga de Expansion'
$ws.Range('D180').Value = 'Mexico Liga de Expansion'
$ws.Range('E180').Value = 45356.92013888889
$ws.Range('F180').Value = 'Atlante'
$ws.Range('G180').Value = 'Tlaxcala FC'
$ws.Range('H180').Value = 2
$ws.Range('I180').Value = 0
$ws.Range('J180').Value = 'H'
$ws.Range('K180').Value = 1.4
$ws.Range('L180').Value = 4.5
$ws.Range('M180').Value = 5.75
$ws.Range('N180').Value = 1.45
$ws.Range('O180').Value = 4.333
$ws.Range('P180').Value = 7.5
$ws.Range('Q180').Value = -1.25
$ws.Range('R180').Value = 1.95
$ws.Range('S180').Value = 1.85
$ws.Range('T180').Value = 2.5
$ws.Range('U180').Value = 1.975
$ws.Range('V180').Value = 1.825
$ws.Range('W180').Value = 0.45
$ws.Range('X180').Value = -1
$ws.Range('Y180').Value = -1
$ws.Range('Z180').Value = 0.95
$ws.Range('AA180').Value = -1
$ws.Range('AB180').Value = -1
$ws.Range('AC180').Value = 0.825

# Row 181
$ws.Range('A181').Value = 179
$ws.Range('B181').Value = 7641687
$ws.Range('C181').Value = 'Mexico Liga de Expansion'
$ws.Range('D181').Value = 'Mexico Liga de Expansion'
$ws.Range('E181').Value = 45357.00347222222
$ws.Range('F181').Value = 'Club Celaya'
$ws.Range('G181').Value = 'Venados FC'
$ws.Range('H181').Value = 0
$ws.Range('I181').Value = 2
$ws.Range('J181').Value = 'A'
$ws.Range('K181').Value = 1.65
$ws.Range('L181').Value = 3.5
$ws.Range('M181').Value = 4.5
$ws.Range('N181').Value = 1.5
$ws.Range('O181').Value = 4
$ws.Range('P181').Value = 7
$ws.Range('Q181').Value = -1
$ws.Range('R181').Value = 1.9
$ws.Range('S181').Value = 1.9
$ws.Range('T181').Value = 2.5
$ws.Range('U181').Value = 1.825
$ws.Range('V181').Value = 1.975
$ws.Range('W181').Value = -1
$ws.Range('X181').Value = -1
$ws.Range('Y181').Value = 6
$ws.Range('Z181').Value = -1
$ws.Range('AA181').Value = 0.8999999999999999
$ws.Range('AB181').Value = -1
$ws.Range('AC181').Value = 0.9750000000000001

# Row 182
$ws.Range('A182').Value = 180
$ws.Range('B182').Value = 7641691
$ws.Range('C182').Value = 'Mexico Liga de Expansion'
$ws.Range('D182').Value = 'Mexico Liga de Expansion'
$ws.Range('E182').Value = 45359.00347222222
$ws.Range('F182').Value = 'Tapatio'
$ws.Range('G182').Value = 'Oaxaca'
$ws.Range('K182').Value = 1.571
$ws.Range('L182').Value = 3.75
$ws.Range('M182').Value = 4.75
$ws.Range('N182').Value = 1.65
$ws.Range('O182').Value = 3.8
$ws.Range('P182').Value = 5
$ws.Range('Q182').Value = -0.75
$ws.Range('R182').Value = 1.775
$ws.Range('S182').Value = 2.025
$ws.Range('T182').Value = 2.5
$ws.Range('U182').Value = 1.85
$ws.Range('V182').Value = 1.95
$ws.Range('W182').Value = 0
$ws.Range('X182').Value = 0
$ws.Range('Y182').Value = 0
$ws.Range('Z182').Value = 0
$ws.Range('AA182').Value = 0

# Row 183
$ws.Range('A183').Value = 181
$ws.Range('B183').Value = 7641692
$ws.Range('C183').Value = 'Mexico Liga de Expansion'
$ws.Range('D183').Value = 'Mexico Liga de Expansion'
$ws.Range('E183').Value = 45361.92013888889
$ws.Range('F183').Value = 'Cimarrones de Sonora FC'
$ws.Range('G183').Value = 'Cancun FC'
$ws.Range('K183').Value = 2.15
$ws.Range('L183').Value = 3.2
$ws.Range('M183').Value = 3.1
$ws.Range('N183').Value = 2.15
$ws.Range('O183').Value = 3.1
$ws.Range('P183').Value = 3.1
$ws.Range('Q183').Value = -0.25
$ws.Range('R183').Value = 1.9
$ws.Range('S183').Value = 1.9
$ws.Range('T183').Value = 2.25
$ws.Range('U183').Value = 1.825
$ws.Range('V183').Value = 1.975
$ws.Range('W183').Value = 0
$ws.Range('X183').Value = 0
$ws.Range('Y183').Value = 0
$ws.Range('Z183').Value = 0
$ws.Range('AA183').Value = 0
